$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "Förändrad" (column C) date value from 2023-09-15 (45184) to
#    2023-09-17 (45186) for every data row (rows 2-115).
for ($r = 2; $r -le 115; $r++) {
    $ws.Range("C" + $r).Value = 45186
}

# 2) Add the case id ("Beteckning", column A) as the friendly-name second
#    argument to every HYPERLINK() formula in columns S, T, V, W, X, Y for
#    the rows that have those links (rows 2-13).
$linkCols = @("S", "T", "V", "W", "X", "Y")
for ($r = 2; $r -le 13; $r++) {
    $id = $ws.Range("A" + $r).Text
    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $r)
        $f = $cell.Formula
        if ($f -and $f.Length -gt 0 -and -not $f.Contains(", ""$id""")) {
            $newF = $f.Substring(0, $f.Length - 1) + ', "' + $id + '")'
            $cell.Formula = $newF
        }
    }
}
